$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 143, pushing the
# previously-existing rows 143:171 down to 144:172 (dimension grows to
# A1:R172). Use a native row insert so formatting/styles shift along
# with the data, then populate the newly-opened row with the new
# observation's values.
$ws.Rows.Item(143).Insert()

$ws.Range("A143").Value = 3
$ws.Range("B143").Value = "Femacal de La Calera"
$ws.Range("C143").Value = "Coquimbo"
$ws.Range("D143").Value = 44505
$ws.Range("E143").Value = 5
$ws.Range("F143").Value = 100112001
$ws.Range("G143").Value = "Berenjena"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 55
$ws.Range("K143").Value = 8000
$ws.Range("L143").Value = 8000
$ws.Range("M143").Value = 8000
$ws.Range("N143").Value = "$/caja 60 unidades"
$ws.Range("O143").Value = "Región de Arica y Parinacota"
$ws.Range("P143").Value = 133
$ws.Range("Q143").Value = 60
$ws.Range("R143").Value = "Hortaliza"
